$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.062.00'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '2.236.23'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '250.90'
$ws.Range('E5').Value = '  +7.92%  '
$ws.Range('D6').Value = '0.625'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('D7').Value = '71.95'
$ws.Range('E7').Value = '  +4.42%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').Value = '  +6.49%  '
$ws.Range('D10').Value = '41.49'
$ws.Range('E10').Value = '  +15.51%  '
$ws.Range('D11').Value = '0.0977'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').Value = '58.49'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '7.27'
$ws.Range('E13').Value = '  +8.25%  '
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '2.567.99'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').Value = '15.08'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').Value = '0.870'
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').Value = '2.230.54'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').Value = '41.932.29'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = '0.0₃0979'
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').Value = '6.25'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = '73.18'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = '236.33'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +9.77%  '
$ws.Range('E25').Value = '  +15.28%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +8.75%  '
$ws.Range('D28').Value = '10.88'
$ws.Range('E28').Value = '  +8.77%  '
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').Value = '171.78'
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('D31').Value = '20.92'
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('D32').Value = '0.123'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '5.66'
$ws.Range('E33').Value = '  +7.79%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.126'
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').Value = '0.0732'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').Value = '4.74'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').Value = '26.26'
$ws.Range('E37').Value = '  +20.40%  '
$ws.Range('E38').Value = '  +9.24%  '
$ws.Range('E39').Value = '  +13.58%  '
$ws.Range('D40').Value = '2.32'
$ws.Range('E40').Value = '  +3.01%  '
$ws.Range('D41').Value = '6.01'
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').Value = '68.41'
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = '11.90'
$ws.Range('E43').Value = '  +18.92%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.209'
$ws.Range('E44').Value = '  +10.79%  '
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '8.85'
$ws.Range('E46').Value = '  -0.50%  '
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('D48').Value = '4.75'
$ws.Range('E48').Value = '  +9.74%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +7.97%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '1.20'
$ws.Range('E51').Value = '  +2.01%  '
